$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$startRow = $used.Row

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $startRow + $i
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value2 = "System, dnasr281@gmail.com"
    }
}
